# FT232H BOM fix: corrected RX and TX LED circuit
# - RX/TX LEDs now use two 470-ohm resistors (R3, R4) instead of a single R5
# - Fuse F1 Digi-Key part number corrected (F4150CT-ND, package 0603 instead of 1206)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Digikey order"
$ws2 = $wb.Worksheets.Item(2)   # "Build BOM"

# ---------------------------------------------------------------
# Sheet "Digikey order": row 9 (470-ohm resistor line)
#   Add Digi-Key part number, bump quantity 1 -> 2, reference R5 -> "R3, R4"
# ---------------------------------------------------------------
$ws1.Range("A9").Value = "P470GCT-ND"
$ws1.Range("E9").Value = 2
$ws1.Range("G9").Value = "R3, R4"

# ---------------------------------------------------------------
# Sheet "Digikey order": row 18 (fuse line)
#   Corrected Digi-Key part number and manufacturer part number/package
# ---------------------------------------------------------------
$ws1.Range("I18").Value = "0603L050YR Littelfuse"
$ws1.Range("A18").Value = "F4150CT-ND"

# ---------------------------------------------------------------
# Sheet "Digikey order": header row (row 1) becomes bold
# ---------------------------------------------------------------
$ws1.Range("A1:L1").Font.Bold = $true

# ---------------------------------------------------------------
# Sheet "Build BOM": rows 18/19 (R4/R5 470-ohm) value column left-aligned
# ---------------------------------------------------------------
$ws2.Range("C18:C19").HorizontalAlignment = -4131

# ---------------------------------------------------------------
# Restore / set selections to match final state
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A18").Select()

$ws2.Activate()
$ws2.Range("C38").Select()

$ws1.Activate()
